$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '37.118.07'
Set-TextValue $ws.Range('E2') '  -1.72%  '
Set-TextValue $ws.Range('D3') '2.024.03'
Set-TextValue $ws.Range('E3') '  -3.06%  '
Set-TextValue $ws.Range('E4') '  -0.19%  '
Set-TextValue $ws.Range('D5') '227.18'
Set-TextValue $ws.Range('E5') '  -2.92%  '
Set-TextValue $ws.Range('D6') '0.609'
Set-TextValue $ws.Range('E6') '  -4.14%  '
Set-TextValue $ws.Range('D8') '55.40'
Set-TextValue $ws.Range('E8') '  -4.88%  '
Set-TextValue $ws.Range('D9') '0.382'
Set-TextValue $ws.Range('E9') '  -2.61%  '
Set-TextValue $ws.Range('D10') '0.0792'
Set-TextValue $ws.Range('E10') '  +1.35%  '
Set-TextValue $ws.Range('E11') '  -3.43%  '
Set-TextValue $ws.Range('D12') '2.319.67'
Set-TextValue $ws.Range('E12') '  -3.13%  '
Set-TextValue $ws.Range('D13') '14.35'
Set-TextValue $ws.Range('E13') '  -5.77%  '
Set-TextValue $ws.Range('D14') '20.51'
Set-TextValue $ws.Range('E14') '  -2.97%  '
Set-TextValue $ws.Range('D15') '0.745'
Set-TextValue $ws.Range('E15') '  -4.50%  '
Set-TextValue $ws.Range('D16') '5.18'
Set-TextValue $ws.Range('E16') '  -3.28%  '
Set-TextValue $ws.Range('D17') '2.044.74'
Set-TextValue $ws.Range('E17') '  -2.15%  '
Set-TextValue $ws.Range('D18') '36.988.02'
Set-TextValue $ws.Range('E18') '  -2.09%  '
Set-TextValue $ws.Range('D19') '6.07'
Set-TextValue $ws.Range('E19') '  -0.40%  '
Set-TextValue $ws.Range('D20') '68.81'
Set-TextValue $ws.Range('E20') '  -3.12%  '
Set-TextValue $ws.Range('D21') '0.0₃0838'
Set-TextValue $ws.Range('E21') '  +0.24%  '
Set-TextValue $ws.Range('D22') '223.32'
Set-TextValue $ws.Range('E22') '  -2.67%  '
Set-TextValue $ws.Range('E23') '  +0.15%  '
Set-TextValue $ws.Range('D24') '2.39'
Set-TextValue $ws.Range('E24') '  +0.56%  '
Set-TextValue $ws.Range('E25') '  -5.68%  '
Set-TextValue $ws.Range('B26') 'Cosmos'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D26') '9.38'
Set-TextValue $ws.Range('E26') '  -3.79%  '
Set-TextValue $ws.Range('B27') 'Monero'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D27') '167.39'
Set-TextValue $ws.Range('E27') '  -2.02%  '
Set-TextValue $ws.Range('D28') '0.126'
Set-TextValue $ws.Range('E28') '  -5.61%  '
Set-TextValue $ws.Range('D29') '18.77'
Set-TextValue $ws.Range('E29') '  -4.08%  '
Set-TextValue $ws.Range('D30') '1.33'
Set-TextValue $ws.Range('E30') '  -4.28%  '
Set-TextValue $ws.Range('E31') '  -4.39%  '
Set-TextValue $ws.Range('D32') '4.49'
Set-TextValue $ws.Range('E32') '  -4.35%  '
Set-TextValue $ws.Range('D33') '0.0608'
Set-TextValue $ws.Range('E33') '  -4.51%  '
Set-TextValue $ws.Range('D34') '4.47'
Set-TextValue $ws.Range('E34') '  -2.90%  '
Set-TextValue $ws.Range('D35') '2.37'
Set-TextValue $ws.Range('E35') '  -5.61%  '
Set-TextValue $ws.Range('E36') '  +0.14%  '
Set-TextValue $ws.Range('D37') '0.999'
Set-TextValue $ws.Range('E37') '  -0.17%  '
Set-TextValue $ws.Range('D38') '3.18'
Set-TextValue $ws.Range('E38') '  -4.51%  '
Set-TextValue $ws.Range('D39') '5.38'
Set-TextValue $ws.Range('D40') '1.506.21'
Set-TextValue $ws.Range('E40') '  +3.77%  '
Set-TextValue $ws.Range('D41') '0.0219'
Set-TextValue $ws.Range('E41') '  -7.19%  '
Set-TextValue $ws.Range('E42') '  -1.97%  '
Set-TextValue $ws.Range('D43') '0.0931'
Set-TextValue $ws.Range('E43') '  -3.24%  '
Set-TextValue $ws.Range('B44') 'Aave'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D44') '95.30'
Set-TextValue $ws.Range('E44') '  -6.09%  '
Set-TextValue $ws.Range('B45') 'InjectiveProtocol'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D45') '16.56'
Set-TextValue $ws.Range('E45') '  -0.93%  '
Set-TextValue $ws.Range('E46') '  -5.47%  '
Set-TextValue $ws.Range('D47') '7.16'
Set-TextValue $ws.Range('E47') '  -0.54%  '
Set-TextValue $ws.Range('E48') '  -4.85%  '
Set-TextValue $ws.Range('E49') '  -1.97%  '
Set-TextValue $ws.Range('D50') '2.207.81'
Set-TextValue $ws.Range('E50') '  -3.06%  '
Set-TextValue $ws.Range('E51') '  -10.70%  '
